$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.009.20"
$ws.Range("E2").Value = "  -0.72%  "

$ws.Range("D3").Value = "1.619.56"
$ws.Range("E3").Value = "  -1.46%  "

$ws.Range("E4").Value = "  -0.05%  "

$ws.Range("D5").Value = "214.11"
$ws.Range("E5").Value = "  -1.42%  "

$ws.Range("D6").Value = "0.517"
$ws.Range("E6").Value = "  +0.62%  "

$ws.Range("E7").Value = "  +0.00%  "

$ws.Range("E8").Value = "  -1.36%  "

$ws.Range("D9").Value = "0.0624"
$ws.Range("E9").Value = "  -0.53%  "

$ws.Range("D10").Value = "20.24"
$ws.Range("E10").Value = "  +1.63%  "

$ws.Range("E11").Value = "  -0.48%  "

$ws.Range("D12").Value = "1.612.90"
$ws.Range("E12").Value = "  -1.57%  "

$ws.Range("E13").Value = "  -0.68%  "

$ws.Range("E14").Value = "  -0.61%  "

$ws.Range("D15").Value = "26.988.17"
$ws.Range("E15").Value = "  -0.76%  "

$ws.Range("D16").Value = "64.29"
$ws.Range("E16").Value = "  -4.62%  "

$ws.Range("E17").Value = "  +0.27%  "

$ws.Range("D18").Value = "215.66"
$ws.Range("E18").Value = "  -1.50%  "

$ws.Range("E19").Value = "  -0.01%  "

$ws.Range("D20").Value = "6.86"
$ws.Range("E20").Value = "  +0.78%  "

$ws.Range("E21").Value = "  -0.75%  "

$ws.Range("E22").Value = "  -5.82%  "

$ws.Range("D23").Value = "8.99"
$ws.Range("E23").Value = "  -1.99%  "

$ws.Range("D24").Value = "147.39"
$ws.Range("E24").Value = "  -0.17%  "

$ws.Range("E25").Value = "  -0.16%  "

$ws.Range("E26").Value = "  -3.77%  "

$ws.Range("E27").Value = "  -1.06%  "

$ws.Range("E28").Value = "  -1.31%  "

$ws.Range("E29").Value = "  -1.11%  "

$ws.Range("E31").Value = "  -1.44%  "

$ws.Range("E32").Value = "  -1.65%  "

$ws.Range("D33").Value = "1.340.03"
$ws.Range("E33").Value = "  +6.23%  "

$ws.Range("E34").Value = "  -0.68%  "

$ws.Range("E35").Value = "  -0.35%  "

$ws.Range("D36").Value = "0.0176"
$ws.Range("E36").Value = "  -1.19%  "

$ws.Range("D37").Value = "0.545"
$ws.Range("E37").Value = "  +0.13%  "

$ws.Range("D38").Value = "0.846"
$ws.Range("E38").Value = "  -0.37%  "

$ws.Range("E39").Value = "  -0.04%  "

$ws.Range("D40").Value = "0.801"
$ws.Range("E40").Value = "  -1.05%  "

$ws.Range("D42").Value = "64.96"
$ws.Range("E42").Value = "  +5.12%  "

$ws.Range("B43").Value = "RocketPoolETH"
$ws.Range("C43").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D43").Value = "1.757.16"
$ws.Range("E43").Value = "  -1.49%  "

$ws.Range("B44").Value = "FraxShare"
$ws.Range("C44").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D44").Value = "5.21"
$ws.Range("E44").Value = "  -2.65%  "

$ws.Range("D45").Value = "90.28"
$ws.Range("E45").Value = "  -1.47%  "

$ws.Range("E46").Value = "  +0.23%  "

$ws.Range("D47").Value = "0.849"
$ws.Range("E47").Value = "  +27.75%  "

$ws.Range("E48").Value = "  -1.16%  "

$ws.Range("E49").Value = "  -0.46%  "

$ws.Range("D50").Value = "0.0994"
$ws.Range("E50").Value = "  +1.90%  "

$ws.Range("E51").Value = "  -1.54%  "
